# feat: add 2022-Q3 data
#
# Before: sheet "总计" (summary) + sheet "2021-Q3" (fund holdings for 2021-Q3).
# After:  sheet "总计" gains a new top data row for 2022-Q3 (the old
#         2021-Q3 row shifts down to row 3); the worksheet that used to hold
#         the "2021-Q3" fund table is repurposed in place to hold brand-new
#         "2022-Q3" fund data, and its old "2021-Q3" content is archived onto
#         a freshly inserted sheet named "2021-Q3" placed right after it.

$wb = $excel.ActiveWorkbook

$total = $wb.Worksheets.Item("总计")
$q3    = $wb.Worksheets.Item("2021-Q3")

# --- 1. Archive the existing "2021-Q3" sheet's content (values + styles)
#        onto a brand new sheet, inserted immediately after the source. ---
$q3.Copy($null, $q3)
$archive = $wb.Worksheets.Item($q3.Index + 1)
$archive.Name = "2021-Q3-tmp"

# --- 2. Repurpose the original sheet in place for the new 2022-Q3 data ---
$q3.Cells.Clear()
$q3.Name = "2022-Q3"

# Give the archived sheet its final name.
$archive.Name = "2021-Q3"

# --- 3. Write the new 2022-Q3 fund table into the repurposed sheet ---
# Pull header formatting from the "总计" sheet's header row so the new
# sheet's header / first-column styling matches (bold + bordered, centered).
$total.Range("B1:D1").Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)
$total.Range("A2").Copy()
$q3.Range("A2").PasteSpecial(-4122)

$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

$q3.Range("A2").Value = 0
$q3.Range("C2").Value = "恒生前海港股通高股息低波动指数"
$q3.Range("H2").Value = 5

# The remaining data cells look numeric ("005702", "0.20", "94.22", ...) but
# must be stored as TEXT (leading zeros / trailing zeros must survive, and
# the source file stores them as strings). Plain `.Value =` lets Excel's
# automatic type-inference coerce them to numbers, so instead build each as
# a formula returning the literal string, then flatten formula -> static
# value via copy / paste-values (keeps the text type without leaving any
# NumberFormat-driven style behind on the cell).
$q3.Range("B2").Formula = '="005702"'
$q3.Range("D2").Formula = '="0.20"'
$q3.Range("E2").Formula = '="94.22"'
$q3.Range("F2").Formula = '="2.55"'
$q3.Range("G2").Formula = '="0.0051"'
$q3.Range("B2:G2").Copy()
$q3.Range("B2:G2").PasteSpecial(-4163)

# --- 4. Update the "总计" (summary) sheet: insert a 2022-Q3 row on top,
#        pushing the existing 2021-Q3 row down to row 3 ---
$total.Range("A2:D2").Copy($total.Range("A3:D3"))
$total.Range("A3").Value = 1

$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 1
$total.Range("D2").Value = 0.01

# Restore the originally-active tab ("总计") as the current selection.
$total.Activate()
